{"js": "// \"RPC Explorer\" -> \"Insight Explorer\" (Slovak menu translation update)\nconst body = context.document.body;\nconst searchResults = body.search(\"RPC Explorer\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"Insight Explorer\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# \"RPC Explorer\" -> \"Insight Explorer\" (Slovak menu translation update)\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"RPC Explorer\"\n$find.Replacement.Text = \"Insight Explorer\"\n$find.Execute(\"RPC Explorer\", $false, $false, $false, $false, $false, $true, 1, $false, \"Insight Explorer\", 2)\n"}
